# Fix bad file path: rename the "1jun23" subfolder to "2023-06-01"
# in every file-path entry on Sheet1 (column A, rows 2-20), and move
# the active selection from A21 to A23 to reflect the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldFolder = "data/gc_data/1jun23/"
$newFolder = "data/gc_data/2023-06-01/"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().Contains($oldFolder)) {
        $cell.Value2 = $val.ToString().Replace($oldFolder, $newFolder)
    }
}

# Update the selected/active cell shown when the sheet is opened
$ws.Range("A23").Select()
